$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.629.80'
$ws.Range("E2").Value = '  -3.87%  '
$ws.Range("D3").Value = '3.348.97'
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.45%  '
$ws.Range("D8").Value = '3.347.77'
$ws.Range("E8").Value = '  -1.93%  '
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.58'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.81%  '
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("E12").Value = '  +1.98%  '
$ws.Range("D13").Value = '3.914.31'
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("E14").Value = '  +0.59%  '
$ws.Range("E15").Value = '  -0.61%  '
$ws.Range("D16").Value = '3.344.49'
$ws.Range("E16").Value = '  -2.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.99'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.81%  '
$ws.Range("D18").Value = '60.694.50'
$ws.Range("E18").Value = '  -3.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.75%  '
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '374.94'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("E23").Value = '  +0.45%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = '3.475.35'
$ws.Range("E25").Value = '  -2.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000116'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +19.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.09'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.35%  '
$ws.Range("E32").Value = '  -0.26%  '
$ws.Range("E33").Value = '  +1.05%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("D35").Value = '3.373.95'
$ws.Range("E35").Value = '  -2.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.14'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.66%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.46'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.46%  '
$ws.Range("E38").Value = '  +3.02%  '
$ws.Range("E39").Value = '  +2.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '162.15'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.22%  '
$ws.Range("E41").Value = '  +2.45%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.750'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.85%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.60'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.30%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.21'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.94'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.87'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +12.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.892'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.48%  '
